$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.523120333333334
$ws.Range("H2").Value = 19.569361
$ws.Range("I2").Value = 0.7039464591847725
$ws.Range("J2").Value = 0.7039464591847725
$ws.Range("M2").Value = 26.34406266666667
$ws.Range("N2").Value = 79.032188
$ws.Range("O2").Value = 0.3168347904490542
$ws.Range("P2").Value = 0.3168347904490542
$ws.Range("Q2").Value = 171.8454908435409
$ws.Range("R2").Value = 1546.609417591868
$ws.Range("S2").Value = 0.2230347288831611
$ws.Range("T2").Value = 0.2230347288831611

# Row 3
$ws.Range("G3").Value = 6.523120333333334
$ws.Range("H3").Value = 19.569361
$ws.Range("I3").Value = 0.7039464591847725
$ws.Range("J3").Value = 0.7039464591847725
$ws.Range("O3").Value = 0.0001452797431229321
$ws.Range("P3").Value = 0.0001452797431229321
$ws.Range("Q3").Value = 0.07879711925322223
$ws.Range("R3").Value = 0.7091740732790001
$ws.Range("S3").Value = 0.0001022691607626613
$ws.Range("T3").Value = 0.0001022691607626613

# Row 4
$ws.Range("G4").Value = 6.523120333333334
$ws.Range("H4").Value = 19.569361
$ws.Range("I4").Value = 0.7039464591847725
$ws.Range("J4").Value = 0.7039464591847725
$ws.Range("O4").Value = 0.6830199298078229
$ws.Range("P4").Value = 0.6830199298078229
$ws.Range("Q4").Value = 370.4577231792967
$ws.Range("R4").Value = 3334.11950861367
$ws.Range("S4").Value = 0.4808094611408488
$ws.Range("T4").Value = 0.4808094611408488

# Row 5
$ws.Range("I5").Value = 0.2740496574363094
$ws.Range("J5").Value = 0.2740496574363095
$ws.Range("M5").Value = 26.34406266666667
$ws.Range("N5").Value = 79.032188
$ws.Range("O5").Value = 0.3168347904490542
$ws.Range("P5").Value = 0.3168347904490542
$ws.Range("Q5").Value = 66.90025538616356
$ws.Range("R5").Value = 602.102298475472
$ws.Range("S5").Value = 0.08682846578646819
$ws.Range("T5").Value = 0.08682846578646819

# Row 6
$ws.Range("I6").Value = 0.2740496574363094
$ws.Range("J6").Value = 0.2740496574363095
$ws.Range("O6").Value = 0.0001452797431229321
$ws.Range("P6").Value = 0.0001452797431229321
$ws.Range("S6").Value = 0.00003981386383527457
$ws.Range("T6").Value = 0.00003981386383527457

# Row 7
$ws.Range("I7").Value = 0.2740496574363094
$ws.Range("J7").Value = 0.2740496574363095
$ws.Range("O7").Value = 0.6830199298078229
$ws.Range("P7").Value = 0.6830199298078229
$ws.Range("S7").Value = 0.187181377786006
$ws.Range("T7").Value = 0.187181377786006

# Row 8
$ws.Range("I8").Value = 0.02200388337891808
$ws.Range("J8").Value = 0.02200388337891808
$ws.Range("M8").Value = 26.34406266666667
$ws.Range("N8").Value = 79.032188
$ws.Range("O8").Value = 0.3168347904490542
$ws.Range("P8").Value = 0.3168347904490542
$ws.Range("Q8").Value = 5.371528033670668
$ws.Range("R8").Value = 48.34375230303601
$ws.Range("S8").Value = 0.006971595779424937
$ws.Range("T8").Value = 0.006971595779424937

# Row 9
$ws.Range("I9").Value = 0.02200388337891808
$ws.Range("J9").Value = 0.02200388337891808
$ws.Range("O9").Value = 0.0001452797431229321
$ws.Range("P9").Value = 0.0001452797431229321
$ws.Range("S9").Value = 0.000003196718524996174
$ws.Range("T9").Value = 0.000003196718524996174

# Row 10
$ws.Range("I10").Value = 0.02200388337891808
$ws.Range("J10").Value = 0.02200388337891808
$ws.Range("O10").Value = 0.6830199298078229
$ws.Range("P10").Value = 0.6830199298078229
$ws.Range("S10").Value = 0.01502909088096815
$ws.Range("T10").Value = 0.01502909088096815
